$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Rename "Image" sheet to "DcImage"
$wb.Worksheets.Item("Image").Name = "DcImage"

# Repurpose the "WorkCreation" sheet into "SchemaCreativeWork" with new headers/data,
# then drop the now-redundant "Work" sheet.
$wsWork = $wb.Worksheets.Item("WorkCreation")
$wsWork.Name = "SchemaCreativeWork"

$wsWork.Cells.Clear()

$wsWork.Range("A1").Value = "@id"
$wsWork.Range("B1").Value = "image"
$wsWork.Range("C1").Value = "name"

$wsWork.Range("A2").Value = "ss-work:test"
$wsWork.Range("B2").Value = "ss-image:test"
$wsWork.Range("C2").Value = "Test work"

$wb.Worksheets.Item("Work").Delete()
